# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held the source filename-derived label
# "5-8-2013-14"; correct it to the actual game date "2014-05-08"
# for every data row (rows 2-31).
#
# A leading apostrophe forces Excel to keep the value as literal text
# instead of auto-converting the ISO-like "2014-05-08" string into a
# date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = "BF"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range("$col$r")
    if ($cell.Text -eq "5-8-2013-14") {
        $cell.Value = "'2014-05-08"
    }
}
